# Update need_to_buy.xlsx data rows 2-15 with refreshed values from R.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# row, A (date serial), B, C, D, E, F
$data = @(
    @(2,  45887, 4535.36362203551, 3883.11865816725, 3912, 4904.619583, 14.1822757971555),
    @(3,  45888, 4535.36362203551, 3889.27259487153, 1944, 4904.619583, 96.4386898265007),
    @(4,  45889, 4535.36362203551, 3941.74220135241, 1944, 4904.619583, 98.6249234298708),
    @(5,  45890, 4535.36362203551, 3955.15644824689, 1944, 4904.619583, 99.1838503838075),
    @(6,  45891, 4535.36362203551, 3955.21842351205, 1944, 4904.619583, 99.1864326865222),
    @(7,  45892, 842.01722651768,  1557.94749249725, 1944, 1638.799688, 17.1137480824819),
    @(8,  45893, 743.085185198094, 1561.6524970481,  1944, 1530.080357, 16.860319535417),
    @(9,  45894, 5476.20428821302, 4316.49098483505, 1944, 5817.607943, 113.078943317585),
    @(10, 45895, 5476.20428821302, 4261.63231553677, 1944, 5817.607943, 110.793165430156),
    @(11, 45896, 5476.20428821302, 4234.376144299,   1944, 5817.607943, 109.657491628583),
    @(12, 45897, 5476.20428821302, 4214.54180171264, 1944, 5817.607943, 108.831060687484),
    @(13, 45898, 5476.20428821302, 4204.10212915779, 1944, 5817.607943, 108.396074331032),
    @(14, 45899, 948.243530910975, 1489.36199351995, 1944, 1809.971087, 16.9620645670407),
    @(15, 45900, 832.128236149678, 1473.0431736709,  1944, 1682.368562, 15.8034791467176)
)

foreach ($entry in $data) {
    $r = $entry[0]
    $ws.Cells.Item($r, 1).Value = $entry[1]
    $ws.Cells.Item($r, 2).Value = $entry[2]
    $ws.Cells.Item($r, 3).Value = $entry[3]
    $ws.Cells.Item($r, 4).Value = $entry[4]
    $ws.Cells.Item($r, 5).Value = $entry[5]
    $ws.Cells.Item($r, 6).Value = $entry[6]
}
